# Weekly "Fruta / hortaliza" update: two new daily price rows are added to
# the top of the data block (rows 204-205), pushing all the existing rows
# down by two (old row 204 -> new row 206, ... old row 277 -> new row 279).
#
# All rows in this block share the same static columns (A,B,C,E,F,G,H,N,O,Q,R);
# only D (Fecha), I (Calidad), J (Volumen), K/L/M (precios) and P (Precio $/Kg)
# differ row to row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 204; this shifts rows 204:277
# down to 206:279 and keeps their original content & formatting intact.
$ws.Range("A204:A205").EntireRow.Insert()

# --- New row 204 -------------------------------------------------------
$ws.Cells.Item(204, 1).Value  = 10
$ws.Cells.Item(204, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(204, 3).Value  = "La Araucanía"
$ws.Cells.Item(204, 4).Value  = 44809
$ws.Cells.Item(204, 5).Value  = 9
$ws.Cells.Item(204, 6).Value  = 100112043
$ws.Cells.Item(204, 7).Value  = "Pepino dulce"
$ws.Cells.Item(204, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(204, 9).Value  = "Primera"
$ws.Cells.Item(204, 10).Value = 240
$ws.Cells.Item(204, 11).Value = 18000
$ws.Cells.Item(204, 12).Value = 19000
$ws.Cells.Item(204, 13).Value = 18417
$ws.Cells.Item(204, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(204, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(204, 16).Value = 1023
$ws.Cells.Item(204, 17).Value = 18
$ws.Cells.Item(204, 18).Value = "Hortaliza"

# --- New row 205 -------------------------------------------------------
$ws.Cells.Item(205, 1).Value  = 10
$ws.Cells.Item(205, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(205, 3).Value  = "La Araucanía"
$ws.Cells.Item(205, 4).Value  = 44809
$ws.Cells.Item(205, 5).Value  = 9
$ws.Cells.Item(205, 6).Value  = 100112043
$ws.Cells.Item(205, 7).Value  = "Pepino dulce"
$ws.Cells.Item(205, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(205, 9).Value  = "Segunda"
$ws.Cells.Item(205, 10).Value = 80
$ws.Cells.Item(205, 11).Value = 16000
$ws.Cells.Item(205, 12).Value = 16000
$ws.Cells.Item(205, 13).Value = 16000
$ws.Cells.Item(205, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(205, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(205, 16).Value = 889
$ws.Cells.Item(205, 17).Value = 18
$ws.Cells.Item(205, 18).Value = "Hortaliza"

# Make sure the Fecha (date) cells use the same date/time display format as
# the rest of column D.
$ws.Cells.Item(204, 4).NumberFormat = $ws.Cells.Item(206, 4).NumberFormat
$ws.Cells.Item(205, 4).NumberFormat = $ws.Cells.Item(206, 4).NumberFormat
